$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "modified" timestamp in B15 to reflect the regeneration time.
$ws.Range("B15").Value = "23/06/2022T14:52:27+00:00"

# Append new vocabulary terms generated from the updated Google Sheet export.
$ws.Range("L24:M25").NumberFormat = "@"
$ws.Range("N24:N25").NumberFormat = "@"

$ws.Range("A24").Value = "vocab:1006"
$ws.Range("B24").Value = "bids"
$ws.Range("D24").Value = "Brain Imaging Data Structure"
$ws.Range("E24").Value = "https://bids.neuroimaging.io/specification.html"
$ws.Range("L24").Value = "2022-06-23"
$ws.Range("M24").Value = "2022-06-23"
$ws.Range("N24").Value = "0000-0001-6361-2571"

$ws.Range("A25").Value = "vocab:1007"
$ws.Range("B25").Value = "power spectral analysis"
$ws.Range("D25").Value = "Computation of the power spectrum of the signal"
$ws.Range("E25").Value = "https://www.sciencedirect.com/topics/engineering/power-spectrum"
$ws.Range("L25").Value = "2022-06-23"
$ws.Range("M25").Value = "2022-06-23"
$ws.Range("N25").Value = "0000-0001-6361-2571"
